# Analysis_Dashboard sheet: add a "Date and Time" row at the top and a
# "Cycle Count of battery" row near the bottom, relabel several metrics
# with their units, fix a couple of values, and append two new speed-bucket
# rows ("70-80 km/h" and "80-90 km/h").
#
# Approach: insert one blank row at the very top (this shifts every
# existing row down by one and carries the existing per-cell styling -
# notably the [hh]:mm:ss number format on the "Total time taken for the
# ride" row - along with it), then (re)write every label/value from row 1
# through row 45 so the sheet ends up exactly as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Insert()

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-12 06:55:31.324000 to 2024-03-12 07:43:09.542000"

$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = 0.03292359953703704

$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 27.70941083333333

$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1390.200271792222

$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 39.497

$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 10.219

$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 25

$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 99

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 28.79475029579625

$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 48.27964325133176

$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 74

$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Custom mode`n95.48%`nEco mode`n3.00%"

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 5458.56792

$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -1767.827968368774

$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 0.01160794777777778

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.0008349768799234199

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.357

$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 2.988

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.3690000000000002

$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 24

$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 39

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 15

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 62

$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 59

$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 57

$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 55

$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 93

$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 44

$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 24

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 20

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 53

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.468598774166667

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.0000001434402615805855

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 44

$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 37.44252152830031

$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 3.072485578128919

$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 3.206253657720926

$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 5.09154752947078

$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 6.993562411169635

$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 8.067887300392943

$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 7.779449878772677

$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 12.79157261098571

$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 15.4293119304406

$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0
